$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 266.66666
$ws.Cells.Item(9, 9).Value = 200
$ws.Cells.Item(9, 11).Value = 200
$ws.Cells.Item(9, 13).Value = -31
$ws.Cells.Item(28, 8).Value = 233.33333
$ws.Cells.Item(28, 9).Value = 170
$ws.Cells.Item(28, 10).Value = 550
$ws.Cells.Item(28, 11).Value = 170
$ws.Cells.Item(28, 12).Value = 550
$ws.Cells.Item(28, 13).Value = 315
$ws.Cells.Item(28, 14).Value = -1520
$ws.Cells.Item(33, 8).Value = 128.33333
$ws.Cells.Item(33, 9).Value = 133.35294
$ws.Cells.Item(33, 11).Value = 133.35294
$ws.Cells.Item(33, 13).Value = 95.64706000000001
$ws.Cells.Item(74, 8).Value = 4399.8
$ws.Cells.Item(74, 9).Value = 3000
$ws.Cells.Item(74, 11).Value = 3000
$ws.Cells.Item(74, 13).Value = -2064
$ws.Cells.Item(77, 8).Value = 4399.8
$ws.Cells.Item(77, 9).Value = 3000
$ws.Cells.Item(77, 11).Value = 15000
$ws.Cells.Item(77, 13).Value = -10320
$ws.Cells.Item(80, 8).Value = 2089.3333
$ws.Cells.Item(80, 9).Value = 3234.2
$ws.Cells.Item(80, 10).Value = 1271.5714
$ws.Cells.Item(80, 11).Value = 9702.599999999999
$ws.Cells.Item(80, 12).Value = 3814.7142
$ws.Cells.Item(80, 13).Value = -8704.599999999999
$ws.Cells.Item(80, 14).Value = -5810.7142
$ws.Cells.Item(83, 8).Value = 2089.3333
$ws.Cells.Item(83, 9).Value = 3234.2
$ws.Cells.Item(83, 10).Value = 1271.5714
$ws.Cells.Item(83, 11).Value = 29107.8
$ws.Cells.Item(83, 12).Value = 11444.1426
$ws.Cells.Item(83, 13).Value = -24115.8
$ws.Cells.Item(83, 14).Value = -21428.1426
$ws.Cells.Item(86, 8).Value = 1359.2
$ws.Cells.Item(86, 9).Value = 933
$ws.Cells.Item(86, 10).Value = 1998.5
$ws.Cells.Item(86, 11).Value = 933
$ws.Cells.Item(86, 12).Value = 1998.5
$ws.Cells.Item(86, 13).Value = 190
$ws.Cells.Item(86, 14).Value = -4244.5
$ws.Cells.Item(89, 8).Value = 1359.2
$ws.Cells.Item(89, 9).Value = 933
$ws.Cells.Item(89, 10).Value = 1998.5
$ws.Cells.Item(89, 11).Value = 4665
$ws.Cells.Item(89, 12).Value = 9992.5
$ws.Cells.Item(89, 13).Value = 951
$ws.Cells.Item(89, 14).Value = -21224.5
$ws.Cells.Item(121, 8).Value = 976.4
$ws.Cells.Item(121, 10).Value = 1494.3334
$ws.Cells.Item(121, 12).Value = 4483.0002
$ws.Cells.Item(121, 14).Value = -7977.0002
$ws.Cells.Item(129, 8).Value = 881.09375
$ws.Cells.Item(129, 9).Value = 748.75
$ws.Cells.Item(129, 10).Value = 900
$ws.Cells.Item(129, 11).Value = 2246.25
$ws.Cells.Item(129, 12).Value = 2700
$ws.Cells.Item(129, 13).Value = 2753.75
$ws.Cells.Item(129, 14).Value = -12700

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3567.8154
$ws.Cells.Item(32, 9).Value = 2924.918
$ws.Cells.Item(32, 11).Value = 2924.918
$ws.Cells.Item(32, 13).Value = -2637.918
$ws.Cells.Item(61, 8).Value = 6603
$ws.Cells.Item(61, 9).Value = 7570.643
$ws.Cells.Item(61, 10).Value = 4909.625
$ws.Cells.Item(61, 11).Value = 7570.643
$ws.Cells.Item(61, 12).Value = 4909.625
$ws.Cells.Item(61, 13).Value = -7358.643
$ws.Cells.Item(61, 14).Value = -5333.625
$ws.Cells.Item(74, 8).Value = 1712.3334
$ws.Cells.Item(74, 9).Value = 449.3889
$ws.Cells.Item(74, 11).Value = 449.3889
$ws.Cells.Item(74, 13).Value = 424.6111
$ws.Cells.Item(77, 8).Value = 1712.3334
$ws.Cells.Item(77, 9).Value = 449.3889
$ws.Cells.Item(77, 11).Value = 2246.9445
$ws.Cells.Item(77, 13).Value = 2121.0555
$ws.Cells.Item(97, 8).Value = 1222.65
$ws.Cells.Item(97, 9).Value = 1144.4706
$ws.Cells.Item(97, 11).Value = 1144.4706
$ws.Cells.Item(97, 13).Value = -648.4706000000001
$ws.Cells.Item(132, 8).Value = 1269.1
$ws.Cells.Item(132, 9).Value = 1142.2972
$ws.Cells.Item(132, 10).Value = 2833
$ws.Cells.Item(132, 11).Value = 3426.8916
$ws.Cells.Item(132, 12).Value = 8499
$ws.Cells.Item(132, 13).Value = -896.8915999999999
$ws.Cells.Item(132, 14).Value = -13559
$ws.Cells.Item(134, 8).Value = 54673.75
$ws.Cells.Item(134, 10).Value = 54673.75
$ws.Cells.Item(134, 12).Value = 54673.75
$ws.Cells.Item(134, 14).Value = -64813.75
$ws.Cells.Item(136, 8).Value = 6603
$ws.Cells.Item(136, 9).Value = 7570.643
$ws.Cells.Item(136, 10).Value = 4909.625
$ws.Cells.Item(136, 11).Value = 22711.929
$ws.Cells.Item(136, 12).Value = 14728.875
$ws.Cells.Item(136, 13).Value = -20161.929
$ws.Cells.Item(136, 14).Value = -19828.875

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(18, 8).Value = 73344.336
$ws.Cells.Item(18, 10).Value = 73344.336
$ws.Cells.Item(18, 12).Value = 73344.336
$ws.Cells.Item(18, 14).Value = -74402.336
$ws.Cells.Item(107, 8).Value = 1020.875
$ws.Cells.Item(107, 9).Value = 833.125
$ws.Cells.Item(107, 10).Value = 1208.625
$ws.Cells.Item(107, 11).Value = 833.125
$ws.Cells.Item(107, 12).Value = 1208.625
$ws.Cells.Item(107, 13).Value = 1086.875
$ws.Cells.Item(107, 14).Value = -5048.625
$ws.Cells.Item(134, 8).Value = 2254.55
$ws.Cells.Item(134, 9).Value = 1965.7693
$ws.Cells.Item(134, 10).Value = 2790.8572
$ws.Cells.Item(134, 11).Value = 5897.3079
$ws.Cells.Item(134, 12).Value = 8372.571599999999
$ws.Cells.Item(134, 13).Value = -3362.3079
$ws.Cells.Item(134, 14).Value = -13442.5716

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1914.0333
$ws.Cells.Item(31, 9).Value = 1592.25
$ws.Cells.Item(31, 10).Value = 2128.5557
$ws.Cells.Item(31, 11).Value = 1592.25
$ws.Cells.Item(31, 12).Value = 2128.5557
$ws.Cells.Item(31, 13).Value = -1297.25
$ws.Cells.Item(31, 14).Value = -2718.5557
$ws.Cells.Item(34, 8).Value = 1914.0333
$ws.Cells.Item(34, 9).Value = 1592.25
$ws.Cells.Item(34, 10).Value = 2128.5557
$ws.Cells.Item(34, 11).Value = 1592.25
$ws.Cells.Item(34, 12).Value = 2128.5557
$ws.Cells.Item(34, 13).Value = -1390.25
$ws.Cells.Item(34, 14).Value = -2532.5557
$ws.Cells.Item(58, 8).Value = 1209388.8
$ws.Cells.Item(58, 9).Value = 2071475
$ws.Cells.Item(58, 11).Value = 2071475
$ws.Cells.Item(58, 13).Value = -2071272
$ws.Cells.Item(94, 8).Value = 1626.3334
$ws.Cells.Item(94, 9).Value = 1796
$ws.Cells.Item(94, 10).Value = 1541.5
$ws.Cells.Item(94, 11).Value = 1796
$ws.Cells.Item(94, 12).Value = 1541.5
$ws.Cells.Item(94, 13).Value = -1345
$ws.Cells.Item(94, 14).Value = -2443.5
$ws.Cells.Item(134, 8).Value = 1465.2122
$ws.Cells.Item(134, 9).Value = 1451.1875
$ws.Cells.Item(134, 11).Value = 4353.5625
$ws.Cells.Item(134, 13).Value = -1818.5625
$ws.Cells.Item(136, 8).Value = 1209388.8
$ws.Cells.Item(136, 9).Value = 2071475
$ws.Cells.Item(136, 11).Value = 6214425
$ws.Cells.Item(136, 13).Value = -6211875

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 142.85715
$ws.Cells.Item(4, 9).Value = 142.85715
$ws.Cells.Item(4, 11).Value = 428.57145
$ws.Cells.Item(4, 13).Value = -316.57145
$ws.Cells.Item(56, 8).Value = 6712.9473
$ws.Cells.Item(56, 9).Value = 6712.9473
$ws.Cells.Item(56, 11).Value = 6712.9473
$ws.Cells.Item(56, 13).Value = -6182.9473
$ws.Cells.Item(69, 8).Value = 2771.9092
$ws.Cells.Item(69, 9).Value = 2166.3333
$ws.Cells.Item(69, 10).Value = 2999
$ws.Cells.Item(69, 11).Value = 6498.999899999999
$ws.Cells.Item(69, 12).Value = 8997
$ws.Cells.Item(69, 13).Value = -5687.999899999999
$ws.Cells.Item(69, 14).Value = -10619
$ws.Cells.Item(72, 8).Value = 2771.9092
$ws.Cells.Item(72, 9).Value = 2166.3333
$ws.Cells.Item(72, 10).Value = 2999
$ws.Cells.Item(72, 11).Value = 19496.9997
$ws.Cells.Item(72, 12).Value = 26991
$ws.Cells.Item(72, 13).Value = -15440.9997
$ws.Cells.Item(72, 14).Value = -35103
$ws.Cells.Item(97, 8).Value = 2800
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 2800
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).ClearContents()
$ws.Cells.Item(97, 13).Value = 8400
$ws.Cells.Item(97, 14).Value = -9392
$ws.Cells.Item(116, 8).Value = 2798.7
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 2798.7
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).ClearContents()
$ws.Cells.Item(116, 13).Value = 8396.099999999999
$ws.Cells.Item(116, 14).Value = -15280.1
$ws.Cells.Item(118, 8).Value = 1037.0834
$ws.Cells.Item(118, 9).Value = 556.875
$ws.Cells.Item(118, 10).Value = 1997.5
$ws.Cells.Item(118, 11).Value = 1670.625
$ws.Cells.Item(118, 12).Value = 5992.5
$ws.Cells.Item(118, 13).Value = -427.625
$ws.Cells.Item(118, 14).Value = -8478.5
$ws.Cells.Item(131, 8).Value = 21333.824
$ws.Cells.Item(131, 10).Value = 24084
$ws.Cells.Item(131, 12).Value = 72252
$ws.Cells.Item(131, 14).Value = -82332

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 655.13336
$ws.Cells.Item(97, 9).Value = 660.2308
$ws.Cells.Item(97, 11).Value = 660.2308
$ws.Cells.Item(97, 13).Value = -164.2308
$ws.Cells.Item(113, 8).Value = 712.7143
$ws.Cells.Item(113, 9).Value = 339.57144
$ws.Cells.Item(113, 11).Value = 339.57144
$ws.Cells.Item(113, 13).Value = 1830.42856
$ws.Cells.Item(116, 8).Value = 40000
$ws.Cells.Item(116, 10).Value = 40000
$ws.Cells.Item(116, 12).Value = 40000
$ws.Cells.Item(116, 14).Value = -49178
$ws.Cells.Item(135, 8).Value = 78988
$ws.Cells.Item(135, 10).Value = 78988
$ws.Cells.Item(135, 12).Value = 78988
$ws.Cells.Item(135, 14).Value = -89128

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 12).ClearContents()
$ws.Cells.Item(14, 14).Value = 0
$ws.Cells.Item(22, 8).Value = 2065.0908
$ws.Cells.Item(22, 9).Value = 1307.8334
$ws.Cells.Item(22, 10).Value = 2973.8
$ws.Cells.Item(22, 11).Value = 1307.8334
$ws.Cells.Item(22, 12).Value = 2973.8
$ws.Cells.Item(22, 13).Value = -1012.8334
$ws.Cells.Item(22, 14).Value = -3563.8
$ws.Cells.Item(27, 8).Value = 2065.0908
$ws.Cells.Item(27, 9).Value = 1307.8334
$ws.Cells.Item(27, 10).Value = 2973.8
$ws.Cells.Item(27, 11).Value = 1307.8334
$ws.Cells.Item(27, 12).Value = 2973.8
$ws.Cells.Item(27, 13).Value = -1200.8334
$ws.Cells.Item(27, 14).Value = -3187.8
$ws.Cells.Item(132, 8).Value = 1606.775
$ws.Cells.Item(132, 9).Value = 1282.56
$ws.Cells.Item(132, 11).Value = 3847.68
$ws.Cells.Item(132, 13).Value = -1317.68
$ws.Cells.Item(136, 8).Value = 2746.6
$ws.Cells.Item(136, 9).Value = 2699.5557
$ws.Cells.Item(136, 11).Value = 8098.6671
$ws.Cells.Item(136, 13).Value = -5548.6671

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 515.13043
$ws.Cells.Item(113, 9).Value = 341.69232
$ws.Cells.Item(113, 11).Value = 1025.07696
$ws.Cells.Item(113, 13).Value = 1144.92304
$ws.Cells.Item(132, 8).Value = 1950.1875
$ws.Cells.Item(132, 9).Value = 1411.2858
$ws.Cells.Item(132, 10).Value = 2979
$ws.Cells.Item(132, 11).Value = 4233.857400000001
$ws.Cells.Item(132, 12).Value = 8937
$ws.Cells.Item(132, 13).Value = -1703.857400000001
$ws.Cells.Item(132, 14).Value = -13997
